# "inventory and shop fixed"
#
# - Row 17 (A17, "tache 19: interface combat") is now marked done (green fill).
# - Row 22 (A22, "mission 5.2: barre d'endu") is now marked done (green fill),
#   and its assignee (B22) changes from "Roméo" to "Fabio" (blue fill).
# - Row 23 (B23, "Bonus: ASCII Art" assignee) changes from "les deux" to
#   "Roméo" (red fill); "les deux" becomes unused and drops out of the
#   shared-string table on save.
# - Selection moves from D14 to D10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reuse the already-existing "done" (green), "Fabio" (blue) and "Roméo" (red)
# cell formats by copy/pasting formats from cells that already carry them,
# instead of assigning raw colors (which would create brand-new style
# entries rather than reusing the workbook's existing ones).

# A17 -> mark done (green), same as A2's format.
$ws.Range("A2").Copy()
$ws.Range("A17").PasteSpecial(-4122)

# A22 -> mark done (green), same as A2's format.
$ws.Range("A2").Copy()
$ws.Range("A22").PasteSpecial(-4122)

# B22 -> reassign to Fabio, using B2's format (Fabio's blue).
$ws.Range("B2").Copy()
$ws.Range("B22").PasteSpecial(-4122)
$ws.Range("B22").Value = "Fabio"

# B23 -> reassign to Roméo, using B3's format (Roméo's red).
$ws.Range("B3").Copy()
$ws.Range("B23").PasteSpecial(-4122)
$ws.Range("B23").Value = "Roméo"

$excel.CutCopyMode = $false

# Update the saved selection.
$ws.Range("D10").Select()
